# Commit: "Fruta / hortaliza, semanal"
#
# The underlying data table (rows 509-573) gets three brand-new weekly
# price records inserted right at its head (new rows 509-511, holding
# "Clementina" / "$/bins (450 kilos)" observations), pushing all of the
# existing records down by three rows (509-573 -> 512-576).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at the top of the data block; Excel shifts
# everything from row 509 down to row 512, carrying formatting
# (including the date number format in column D) along with it.
$ws.Rows("509:511").Insert()

# Columns that are constant across every "Fruta - Mandarina" record in
# this sheet (Mercado ID, Mercado, Región, Codreg, Tipo, Producto ID,
# Producto, Categoría ID, Categoría, Origen).
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100102
$producto   = "Cítricos"
$categoriaId = 100102004
$categoria   = "Mandarina"
$origen      = "Provincia de Limarí"

function Set-Row($r, $fecha, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $unidad, $precioKg, $kgUnidad) {
    $ws.Range("A$r").Value = $mercadoId
    $ws.Range("B$r").Value = $mercado
    $ws.Range("C$r").Value = $region
    $ws.Range("D$r").Value = $fecha
    $ws.Range("E$r").Value = $codreg
    $ws.Range("F$r").Value = $tipo
    $ws.Range("G$r").Value = $productoId
    $ws.Range("H$r").Value = $producto
    $ws.Range("I$r").Value = $categoriaId
    $ws.Range("J$r").Value = $categoria
    $ws.Range("K$r").Value = $variedad
    $ws.Range("L$r").Value = $calidad
    $ws.Range("M$r").Value = $volumen
    $ws.Range("N$r").Value = $precioMin
    $ws.Range("O$r").Value = $precioMax
    $ws.Range("P$r").Value = $precioProm
    $ws.Range("Q$r").Value = $unidad
    $ws.Range("R$r").Value = $origen
    $ws.Range("S$r").Value = $precioKg
    $ws.Range("T$r").Value = $kgUnidad
}

# New row 509: Clementina / Especial, $/bins (450 kilos), fecha 45077 (2023-05-31)
Set-Row 509 45077 "Clementina" "Especial" 20  260000 270000 265000 "$/bins (450 kilos)" 589 450

# New row 510: Clementina / Primera
Set-Row 510 45077 "Clementina" "Primera"  24  220000 230000 225000 "$/bins (450 kilos)" 500 450

# New row 511: Clementina / Segunda
Set-Row 511 45077 "Clementina" "Segunda"  20  190000 200000 195000 "$/bins (450 kilos)" 433 450
